$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row at the top of the data table (row 6) for the new
# "Febrero 2025" entry, pushing all existing rows down by one.
$ws.Rows.Item(6).Insert()

# Resize the table / autofilter to include the new row.
$lo.Resize($ws.Range("B5:F91"))

# Populate the new row with the Feb-2025 figures.
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Feb."
$ws.Range("D6").Formula = "=SUM(E6:F6)"
$ws.Range("E6").Value = 4481943
$ws.Range("F6").Value = 4884035

# Clear out the Sep-2022 figures (row 35 after the insert) - the source
# data for that month was retracted.
$ws.Range("E35").ClearContents()
$ws.Range("F35").ClearContents()

# Update the "last updated" caption.
$ws.Range("B92").Value = "Actualización: Febrero 2025."
